$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header F1: "Number of Borrowed Books" -> "Borrowed Books"
$ws.Range("F1").Value = "Borrowed Books"

# Insert two new student rows right after the current row 3 (a246810),
# before the old row 4 (a55664478) - they become new rows 4 and 5,
# pushing the existing rows 4-7 down to 6-9.
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# New row 4 data
$ws.Range("A4").Value = "a32132111231"
$ws.Range("B4").Value = "sadadasd"
$ws.Range("C4").Value = "HSU"
$ws.Range("D4").Value = "sadad"
# Force the all-digit contact number to be stored as text (not a number),
# without leaving behind a new (unused) number-format style: compute it as
# text via a formula, then paste-special as a value over itself.
$ws.Range("E4").Formula = '=TEXT(321313123313,"0")'
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)
$ws.Range("F4").Value = 1

# New row 5 data
$ws.Range("A5").Value = "a342423244"
$ws.Range("B5").Value = "sdaasdsadasdsa"
$ws.Range("C5").Value = "I"
$ws.Range("D5").Value = "tcsada"
# Force the all-digit contact number to be stored as text (not a number)
$ws.Range("E5").Formula = '=TEXT(321313131212,"0")'
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)
$ws.Range("F5").Value = 1

# Set "Borrowed Books" count to 1 for every remaining existing student row
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = 1
